$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Re-point every "old" table style onto the new table style GUID.
#    (6 tables across the deck share the same old style id.)
# ---------------------------------------------------------------------------
$oldStyleId = "{22567273-7B7A-4A14-B2B7-553CC2E4EEBA}"
$newStyleId = "{50340CCA-2B8B-4F2A-B127-A3B99710AB17}"

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $s = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTable) {
            $tbl = $sh.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Swap the "Default" / "Simple Light" colour schemes between the two
#    auxiliary themes. The COM surface only exposes the presentation's
#    primary theme colours through Slide.ThemeColorScheme, which is backed
#    by theme2.xml ("Simple Light" -> becomes "Default").
# ---------------------------------------------------------------------------
function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$tcs = $p.Slides.Item(1).ThemeColorScheme

# dk1/lt1 (black/white) are unchanged by the edit.
$tcs.Colors(3).RGB  = RGBVal 0x15 0x81 0x58   # dk2      -> 158158
$tcs.Colors(4).RGB  = RGBVal 0xF3 0xF3 0xF3   # lt2      -> F3F3F3
$tcs.Colors(5).RGB  = RGBVal 0x05 0x8D 0xC7   # accent1  -> 058DC7
$tcs.Colors(6).RGB  = RGBVal 0x50 0xB4 0x32   # accent2  -> 50B432
$tcs.Colors(7).RGB  = RGBVal 0xED 0x56 0x1B   # accent3  -> ED561B
$tcs.Colors(8).RGB  = RGBVal 0xED 0xEF 0x00   # accent4  -> EDEF00
$tcs.Colors(9).RGB  = RGBVal 0x24 0xCB 0xE5   # accent5  -> 24CBE5
$tcs.Colors(10).RGB = RGBVal 0x64 0xE5 0x72   # accent6  -> 64E572
$tcs.Colors(11).RGB = RGBVal 0x22 0x00 0xCC   # hlink    -> 2200CC
$tcs.Colors(12).RGB = RGBVal 0x55 0x1A 0x8B   # folHlink -> 551A8B
